# Daily attendance processing - reorder "Recorded By" (column G) entries so
# that any "System" token is moved to the front of the comma-separated list,
# preserving the relative order of the remaining tokens.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$rowCount = $used.Rows.Count

$changedCount = 0

for ($r = 2; $r -le $rowCount; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # Column G = "Recorded By"
    $val = $cell.Value2

    if ($val -eq $null) {
        continue
    }
    if (-not ($val -is [string])) {
        continue
    }
    if ($val -eq "") {
        continue
    }
    if ($val.IndexOf(",") -lt 0) {
        continue
    }

    $parts = $val -split ","
    $trimmed = @()
    foreach ($p in $parts) {
        $trimmed += $p.Trim()
    }

    $systemParts = @()
    $restParts = @()
    foreach ($p in $trimmed) {
        if ($p.Equals("System")) {
            $systemParts += $p
        } else {
            $restParts += $p
        }
    }

    if ($systemParts.Count -eq 0) {
        continue
    }

    $newOrder = $systemParts + $restParts
    $newVal = [string]::Join(", ", $newOrder)

    if (-not $newVal.Equals($val)) {
        $cell.Value2 = $newVal
        $changedCount += 1
    }
}

Write-Output ("Recorded By values reordered: " + $changedCount)
